# Rename the "Requested quantity" headers on the two existing sheets and
# add a new "PO Forecast" sheet with a ds / PO_Forecast / yhat_lower /
# yhat_upper forecast table (Prophet-style output).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Weekly Quantity
$ws2 = $wb.Worksheets.Item(2)   # Monthly Trend

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet right after "Monthly Trend" so it lands last / gets
# sheetId 3, matching the workbook order in the diff.
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "PO Forecast"

$headers = @("ds", "PO_Forecast", "yhat_lower", "yhat_upper")
for ($j = 0; $j -lt $headers.Length; $j++) {
    $ws3.Cells.Item(1, $j + 1).Value = $headers[$j]
}

$data = @(
    @(44983.99999999999, 46, -109.4961925014655, 192.5948956430826),
    @(45004.99999999999, 48, -103.6663804414446, 199.6607635808974),
    @(45011.99999999999, 49, -109.1537271305042, 192.5772926906727),
    @(45081.99999999999, 58, -94.11207480802713, 207.8011009399883),
    @(45095.99999999999, 60, -90.17541006055507, 214.67708631043),
    @(45109.99999999999, 62, -103.8235741746622, 208.7936784547367),
    @(45137.99999999999, 65, -96.83460922977427, 210.0893750875715),
    @(45165.99999999999, 69, -73.37369674692638, 214.203221174295),
    @(45172.99999999999, 70, -80.35998655829863, 215.229997959247),
    @(45179.99999999999, 71, -76.72081657738651, 219.4176489319345),
    @(45207.99999999999, 74, -80.29833753445313, 221.0422898039483),
    @(45557.99999999999, 119, -38.22937753294465, 277.0993599454376),
    @(45564.99999999999, 120, -26.19323014548186, 272.9999222107436),
    @(45571.99999999999, 121, -35.67267191236229, 275.7934634630042),
    @(45578.99999999999, 122, -35.33961253536719, 269.9457848205878),
    @(45585.99999999999, 123, -22.50548825145068, 269.0882871658732),
    @(45592.99999999999, 123, -12.1148871033229, 287.6350718554778),
    @(45599.99999999999, 124, -22.84122371434833, 279.4213573762946),
    @(45606.99999999999, 125, -27.34256790557204, 262.9041691528553),
    @(45613.99999999999, 126, -37.40252222232969, 265.6803998962648)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 1).Value = $data[$i][0]
    $ws3.Cells.Item($row, 2).Value = $data[$i][1]
    $ws3.Cells.Item($row, 3).Value = $data[$i][2]
    $ws3.Cells.Item($row, 4).Value = $data[$i][3]
}

# Reuse the existing header / date styles instead of minting new ones, to
# mirror the "s=1" / "s=2" style reuse seen in the target workbook.
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws3.Range("A2:A21").PasteSpecial(-4122)

# Restore the original active sheet (the workbook opened on "Weekly
# Quantity"); adding the new sheet otherwise leaves it focused last.
$ws1.Activate() | Out-Null
